$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B2" = 27.74653094557529
    "C2" = 8.002295477005312
    "D2" = 13.14626756566572
    "E2" = 12.94890160155014
    "G2" = 65.94053527659115
    "H2" = 23.36465109646661
    "J2" = 7.970545630622736
    "L2" = 13.70671795778608
    "N2" = 21.32934096539442
    "B3" = 27.44392620247561
    "C3" = 7.61719194961475
    "D3" = 13.15658724123946
    "E3" = 12.9780167384613
    "G3" = 65.76225826713453
    "H3" = 23.39017951330844
    "J3" = 7.973165359565686
    "L3" = 13.70812507771486
    "N3" = 21.40356120827893
    "B4" = 27.26384658482954
    "C4" = 7.37261184371141
    "D4" = 13.16525635019611
    "E4" = 12.99726767434149
    "G4" = 65.67132567530301
    "H4" = 23.41071208423253
    "J4" = 7.974898548857663
    "L4" = 13.71138630595154
    "N4" = 21.45117735915216
    "B5" = 27.19198097252418
    "C5" = 7.271048406974987
    "D5" = 13.16937526748034
    "E5" = 13.00545854883766
    "G5" = 65.63894034971034
    "H5" = 23.42029752347036
    "J5" = 7.975636246462995
    "L5" = 13.71331846051798
    "N5" = 21.47109690267346
    "B6" = 27.18014170934904
    "C6" = 7.254074180271308
    "D6" = 13.17009459943566
    "E6" = 13.00683954517147
    "G6" = 65.63384508337516
    "H6" = 23.42196265278312
    "J6" = 7.975760639407266
    "L6" = 13.71367573011215
    "N6" = 21.47443570952782
    "B7" = 27.26287112973009
    "C7" = 7.371249584804701
    "D7" = 13.16530952658044
    "E7" = 12.99737673798646
    "G7" = 65.67086999508754
    "H7" = 23.41083642899324
    "J7" = 7.974908370450326
    "L7" = 13.71140992111269
    "N7" = 21.45144391169907
    "B8" = 27.64105423300599
    "C8" = 7.871280280087231
    "D8" = 13.14934141596401
    "E8" = 12.95865555760059
    "G8" = 65.87522758148752
    "H8" = 23.37244310839958
    "J8" = 7.971423075312842
    "L8" = 13.70670586724156
    "N8" = 21.35450861433416
    "B9" = 28.42441975929746
    "C9" = 8.781918114307134
    "D9" = 13.13655110881492
    "E9" = 12.89360761174305
    "G9" = 66.42227887300666
    "H9" = 23.33583625211064
    "J9" = 7.965574912288269
    "L9" = 13.71647875114375
    "N9" = 21.18057477447626
    "B10" = 29.02038063976316
    "C10" = 9.402327191499715
    "D10" = 13.13845735341362
    "E10" = 12.85242660307062
    "G10" = 66.91207257172829
    "H10" = 23.33268954078639
    "J10" = 7.961876068862623
    "L10" = 13.73519772092978
    "N10" = 21.06253854997607
    "B11" = 29.29483868187673
    "C11" = 9.673053638776491
    "D11" = 13.14177784139085
    "E11" = 12.83512222302519
    "G11" = 67.15358034117219
    "H11" = 23.33644217542937
    "J11" = 7.960322431834283
    "L11" = 13.74620515406463
    "N11" = 21.01093811371721
    "B12" = 29.39915882953883
    "C12" = 9.773854751022446
    "D12" = 13.14338754422771
    "E12" = 12.82877460519675
    "G12" = 67.24768015384464
    "H12" = 23.33861010104243
    "J12" = 7.959752600621929
    "L12" = 13.75073009710528
    "N12" = 20.99169812432258
    "B13" = 29.37667580245688
    "C13" = 9.752222689789306
    "D13" = 13.14302520668762
    "E13" = 12.83013255812336
    "G13" = 67.22729712395193
    "H13" = 23.33810996379954
    "J13" = 7.959874502086792
    "L13" = 13.74973973638275
    "N13" = 20.9958284786021
    "B14" = 29.30341390186546
    "C14" = 9.681381344375948
    "D14" = 13.14190321529997
    "E14" = 12.8345958904965
    "G14" = 67.16126922068224
    "H14" = 23.33660555679596
    "J14" = 7.960275181047704
    "L14" = 13.74657028628853
    "N14" = 21.00934922339468
    "B15" = 29.25858675518068
    "C15" = 9.63776366376314
    "D15" = 13.14126182633456
    "E15" = 12.83735652084926
    "G15" = 67.12116846617172
    "H15" = 23.33578136479271
    "J15" = 7.960523015906183
    "L15" = 13.7446753007069
    "N15" = 21.01767009865149
    "B16" = 29.00250363836448
    "C16" = 9.384397488254054
    "D16" = 13.13828968233036
    "E16" = 12.85358621377186
    "G16" = 66.89666218927519
    "H16" = 23.33254874600837
    "J16" = 7.96198019387107
    "L16" = 13.73452836746568
    "N16" = 21.06595280524066
    "B17" = 28.84619709148043
    "C17" = 9.22597126438659
    "D17" = 13.13709444392535
    "E17" = 12.86390839169166
    "G17" = 66.76369512015043
    "H17" = 23.33189456130786
    "J17" = 7.962907125358981
    "L17" = 13.72894071032449
    "N17" = 21.09610831769949
    "B18" = 28.75661481273151
    "C18" = 9.133769199183844
    "D18" = 13.13663793408678
    "E18" = 12.86997996499482
    "G18" = 66.68897882210129
    "H18" = 23.33200618643177
    "J18" = 7.963452415516316
    "L18" = 13.72596152125121
    "N18" = 21.11365022078306
    "B19" = 28.72634179314262
    "C19" = 9.102367874691748
    "D19" = 13.13652304670253
    "E19" = 12.87205881088711
    "G19" = 66.66398515611779
    "H19" = 23.33212772647578
    "J19" = 7.963639128637766
    "L19" = 13.72499317354435
    "N19" = 21.11962352421936
    "B20" = 28.86280363351629
    "C20" = 9.242948223597844
    "D20" = 13.13719777851462
    "E20" = 12.8627956572904
    "G20" = 66.7776675375993
    "H20" = 23.3319136934236
    "J20" = 7.962807195464062
    "L20" = 13.72951124857305
    "N20" = 21.09287780648778
    "B21" = 29.32492288456715
    "C21" = 9.702236217262614
    "D21" = 13.14222321465297
    "E21" = 12.83327933573131
    "G21" = 67.18059176804533
    "H21" = 23.33702715831145
    "J21" = 7.960156990307812
    "L21" = 13.74749156548523
    "N21" = 21.00536972157607
    "B22" = 29.62916835455729
    "C22" = 9.992376827148378
    "D22" = 13.14756077938177
    "E22" = 12.81518448627942
    "G22" = 67.45932605804943
    "H22" = 23.3447228528535
    "J22" = 7.958532722326947
    "L22" = 13.7613207097577
    "N22" = 20.94992604773273
    "B23" = 29.46661394259841
    "C23" = 9.838459289146556
    "D23" = 13.1445243650333
    "E23" = 12.82473273988834
    "G23" = 67.30916659951743
    "H23" = 23.34021679623947
    "J23" = 7.959389777983031
    "L23" = 13.75375032412345
    "N23" = 20.97935785628925
    "B24" = 28.85529493528514
    "C24" = 9.235276426209756
    "D24" = 13.13715034249356
    "E24" = 12.86329829709248
    "G24" = 66.77134522140588
    "H24" = 23.33190352475351
    "J24" = 7.962852335189654
    "L24" = 13.72925258164807
    "N24" = 21.09433768167896
    "B25" = 28.20854231779372
    "C25" = 8.543737374036741
    "D25" = 13.13802634327413
    "E25" = 12.91004224277124
    "G25" = 66.25873811184711
    "H25" = 23.34157933998805
    "J25" = 7.967051759061329
    "L25" = 13.71180454253272
    "N25" = 21.22590887942397
}

foreach ($key in $updates.Keys) {
    $ws.Range($key).Value = $updates[$key]
}
